$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.895.66'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '2.500.63'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'535.48"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = "'137.96"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("D9").Value = '2.524.82'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("D14").Value = '2.961.12'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = "'23.27"
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").Value = '58.893.84'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '2.508.71'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").Value = "'11.11"
$ws.Range("E19").Value = '  +1.68%  '
$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = "'325.31"
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").Value = "'5.89"
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = "'64.88"
$ws.Range("E24").Value = '  +4.38%  '
$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("D28").Value = "'7.63"
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("D29").Value = '0.0₃0776'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").Value = "'6.72"
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = "'1.76"
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").Value = "'168.12"
$ws.Range("E32").Value = '  +3.62%  '
$ws.Range("E33").Value = '  +5.69%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").Value = "'1.41"
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").Value = "'18.61"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").Value = "'4.11"
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("D39").Value = "'36.74"
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("E40").Value = '  +3.83%  '
$ws.Range("D41").Value = "'3.64"
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").Value = "'5.30"
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = "'282.25"
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Value = "'0.995"
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").Value = "'130.94"
$ws.Range("E45").Value = '  +6.93%  '
$ws.Range("D46").Value = "'0.606"
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = "'0.0512"
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = "'17.40"
$ws.Range("E51").Value = '  -0.14%  '
